$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "2GJWGR"
$ws.Range("B5").Value = "2024-11-18 09:30:32"
$ws.Range("C5").Value = "GET /"
$ws.Range("D5").Value = 401
$ws.Range("E5").Value = $false
$f5 = @"
Eitss... mau ngapain? Akses terbatas!
"@
$ws.Range("F5").Value = $f5

$ws.Range("A6").Value = "A870XX"
$ws.Range("B6").Value = "2024-11-18 09:30:32"
$ws.Range("C6").Value = "GET /favicon.ico"
$ws.Range("D6").Value = 404
$ws.Range("E6").Value = $false
$f6 = @"
Not Found
"@
$ws.Range("F6").Value = $f6

$ws.Range("A7").Value = "IKS2M5"
$ws.Range("B7").Value = "2024-11-18 09:31:10"
$ws.Range("C7").Value = "GET /"
$ws.Range("D7").Value = 200
$ws.Range("E7").Value = $true
$f7 = @"
API Virtual Assistant Undiksha 2024-11-18 09:31:10
"@
$ws.Range("F7").Value = $f7

$ws.Range("A8").Value = "GZ455E"
$ws.Range("B8").Value = "2024-11-18 09:31:14"
$ws.Range("C8").Value = "GET /"
$ws.Range("D8").Value = 200
$ws.Range("E8").Value = $true
$f8 = @"
API Virtual Assistant Undiksha 2024-11-18 09:31:14
"@
$ws.Range("F8").Value = $f8

$ws.Range("A9").Value = "N0RUG7"
$ws.Range("B9").Value = "2024-11-18 09:32:22"
$ws.Range("C9").Value = "POST /chat"
$ws.Range("D9").Value = 200
$ws.Range("E9").Value = $true
$f9 = @"
OK
###
Question:Siapa rektor undiksha? Saya ingin cetak ktm 2115101014. Saya ingin cek kelulusan nomor pendaftaran 3242000006 tanggal lahir 2005-11-30. Saya ingin bunuh diri.
###
Answer:Salam Harmoni🙏
Berikut adalah informasi yang dapat disampaikan:
1. Pertanyaan yang diajukan tidak relevan dengan konteks kampus Universitas Pendidikan Ganesha.
2. Informasi mengenai kelulusan peserta SMBJM di Universitas Pendidikan Ganesha:
   - Nomor Pendaftaran: 3242000006
   - Nama Siswa: KADEK YUNI SEDANI
   - Tanggal Lahir: 30 November 2005
   - Tahun Daftar: 2024
   - Pilihan Program Studi: Pendidikan Jasmani, Kesehatan Dan Rekreasi (S1)
   - Status Kelulusan: Lulus
   Jika KADEK YUNI SEDANI lulus, selamat bergabung menjadi bagian dari Universitas Pendidikan Ganesha! Jika tidak lulus, tetap semangat dan terus berusaha!
3. Informasi mengenai Kartu Tanda Mahasiswa (KTM):
   - NIM: 2115101014
   - URL KTM: [KTM Anda](https://aka.undiksha.ac.id/api/ktm/generate/2115101014?token=VGhoSWZmU3MzSE5HRjJjV0hoazRUNWlGZnUrczdNMWpqUFBFNUpZdDZwbUliSitYN1Vvbm8vNjBBYndlTVpGOTBESU9Kamc3cEV2bW1SYitSbnYxd2U4dkxIVkExczhwMFovb2lOODlJbjdrTkZZSG5XZkMxbWRmVWRXWEE4ZWxYM2lNemtJWnUvRHZ5dCs1VnFCVFBwdndxQU1IZWNkcjVqNU1Jamw4bVdKcFhZSHhxaFM0NnhrdXpKWHhWRHZsVUR4aVFSNWVVaHAvK3FGZGsrYWFnR1N0czY1QUdBMlB6WkUxOXU0Y0VjSm5SZlU0aUlhTWJFTVBBMVVnWGRzU1ZjZXdreTR3YVFLbDVKczU4R3J4T3JVNkE2dDRES25LY3Y0bHBuREtWYUU1Nk9vN0VmMmpkMWJHd2ljYll5bHBWTXQ3Y1Y3ZEVrY0xETTBqZlNOdDhCcWgvdnRHMWh1azBtdUI5VGFWQ1ZMSzNJRjJJYU9DUUlreG9pTG1FSE5iRGs4WHdBbEtqV1pNeE9IVXRLSWJrWDA3RGh3UXRiTmZORGc3dXRPdmhGYjdjM2pqRDVNazdUSGxBelU9).
4. Rektor Universitas Pendidikan Ganesha adalah Prof. Dr. I Wayan Lasmawan, M.Pd. 
Semoga informasi ini bermanfaat!
"@
$ws.Range("F9").Value = $f9

$ws.Range("A10").Value = "1XK6U0"
$ws.Range("B10").Value = "2024-11-18 09:32:53"
$ws.Range("C10").Value = "POST /chat"
$ws.Range("D10").Value = 200
$ws.Range("E10").Value = $true
$f10 = @"
OK
###
Question:Siapa rektor undiksha? Saya ingin cetak ktm 2115101014. Saya ingin cek kelulusan nomor pendaftaran 3242000006 tanggal lahir 2005-11-30. Saya ingin bunuh diri.
###
Answer:Salam Harmoni🙏
Berikut adalah informasi yang dapat disampaikan:
1. Pertanyaan yang diajukan tidak relevan dengan konteks kampus Universitas Pendidikan Ganesha.
2. Informasi mengenai kelulusan peserta SMBJM di Universitas Pendidikan Ganesha:
   - Nomor Pendaftaran: 3242000006
   - Nama Siswa: KADEK YUNI SEDANI
   - Tanggal Lahir: 30 November 2005
   - Tahun Daftar: 2024
   - Pilihan Program Studi: Pendidikan Jasmani, Kesehatan Dan Rekreasi (S1)
   - Status Kelulusan: Lulus
   Jika KADEK YUNI SEDANI lulus, selamat bergabung menjadi bagian dari Universitas Pendidikan Ganesha! Jika tidak lulus, tetap semangat dan terus berusaha!
3. Informasi mengenai Kartu Tanda Mahasiswa (KTM):
   - NIM: 2115101014
   - URL KTM: [KTM Anda](https://aka.undiksha.ac.id/api/ktm/generate/2115101014?token=N291VWNLZVpxSXpWY0pCanVYTUNLQ0hJaGxQdXFZelR6YmhDbENMMTBBbzFNRnpRaXY4M2Zwb0ZuUTVHMk9DRStaY09GVFlISkdnSWxXdy9pOHROWEdWTmpOWit2cEhoWjlaaE1wc3lIRkJ3L3J0dTg1Vk5aRWlrdUtZczQxZFI4VDU5U0VFcy9OTU5qbUlhbDF4OWZZN1NsU3JUY3NvejhXWTBEamxtZ1NPMWNUNy9ZUm9jT2tjRkRQLzE0dExmV29Lc3htemxFZ2trRkpzMnpMTmJGWGdtSm9yRjZEVGttNTZ2bGZ5SHZCSnA4SGFhSzBKK25LSlI5VEo3bEx3T0lNaWp6VytsaEpIdFZzYXlnMGlUVGdMWitPSEFtSmRCelRNd2JoQTdiUGlSN3pZSUR5NkU0a0hzYkxrQW83eGJVeE0va1RkYTlBNHB1NXNZd0cxc3BpdGVDYmxHcVRzeXNIZ3MzUHFHL3NMb0pHQndRbE5LNWQ5TmVydW1BZUNzOG1WUk8xUjNFSERPUU9ENEVJN0JyTDJiRGd2cmdIVXdNREJSaS9UTUp2cXJVTUtYc3hQMXh4bmNjWEE9)
4. Rektor Universitas Pendidikan Ganesha adalah Prof. Dr. I Wayan Lasmawan, M.Pd. 
Semoga informasi ini bermanfaat!
"@
$ws.Range("F10").Value = $f10
